# Leave Card update — "Upload Leave Card 12/27/2023 4:01 PM"
#
# Fills in the monthly "EARNED" (C) and running-balance (E / I) columns
# through November 2023, records a Sick-Leave remark for 11/28-29/2023,
# and lets the table's calculated BALANCE cells (E9 / I9) recompute.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in EARNED (column C) for Apr-2023 .. Nov-2023 rows that were blank ---
$ws.Range("C27").Value = 1.25   # 45017 (2023-04-01)
$ws.Range("C28").Value = 1.25   # 45047 (2023-05-01)
$ws.Range("C29").Value = 1.25   # 45078 (2023-06-01)
$ws.Range("C30").Value = 1.25   # 45108 (2023-07-01)
$ws.Range("C31").Value = 1.25   # 45139 (2023-08-01)
$ws.Range("C32").Value = 1.25   # 45170 (2023-09-01)
$ws.Range("C33").Value = 1.25   # 45200 (2023-10-01)
$ws.Range("C34").Value = 1.25   # 45231 (2023-11-01)

# --- New Sick Leave entry for November 2023 ---
$ws.Range("B34").Value = "SL(2-0-0)"
$ws.Range("H34").Value = 2
$ws.Range("K34").Value = "11/28,29/2023"

# --- Running balance formulas, columns E and I, rows 11-34 ---
# (row 10 already holds the brought-forward balance; these were blank)
$ws.Range("E11").Formula = "=SUM(C11,E10)-D11"
$ws.Range("E12:E54").Formula = "=SUM(C12,E11)-D12"
$ws.Range("E35:E54").ClearContents()

$ws.Range("I11").Formula = "=SUM(G11,I10)-H11"
$ws.Range("I12:I75").Formula = "=SUM(G12,I11)-H12"
$ws.Range("I35:I75").ClearContents()

$excel.Calculate()
